$wb = $excel.ActiveWorkbook

# --- Overview sheet: mark a.md as handed back (in sync with en-US) for both locales ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: record the new handback for a.md (row 2) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("L2").Value = "2017-02-21 09:49:09"
$zhcn.Range("M2").Value = "TestHandback_201702210549"
$zhcn.Range("R2").Value = ""

# --- de-de sheet: record the new handback for a.md (row 2) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("L2").Value = "2017-02-21 09:49:31"
$dede.Range("M2").Value = "TestHandback_201702210549"
$dede.Range("R2").Value = ""

# --- Autofit the columns that now contain longer text ---
$overview.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$overview.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$zhcn.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$zhcn.Columns.Item(13).EntireColumn.AutoFit() | Out-Null
$dede.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$dede.Columns.Item(13).EntireColumn.AutoFit() | Out-Null
